# Generate Report for Handback
# - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   on the Overview sheet (both locale columns) and on each locale sheet.
# - Each locale sheet's "Latest Handback DateTime" is refreshed to the
#   handback run time.
# - The stale "handback file is not latest" error is cleared now that the
#   handback is in sync.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet --------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E2").Value = $newStatus
$ws.Range("F2").Value = $newStatus
$ws.Columns.Item(5).AutoFit() | Out-Null
$ws.Columns.Item(6).AutoFit() | Out-Null

# --- zh-cn sheet -------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("K2").Value = "2016-08-21 14:56:25"
$wsZh.Range("P2").Value = ""
$wsZh.Columns.Item(3).AutoFit() | Out-Null
$wsZh.Columns.Item(16).AutoFit() | Out-Null

# --- de-de sheet -------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("K2").Value = "2016-08-21 14:56:34"
$wsDe.Range("P2").Value = ""
$wsDe.Columns.Item(3).AutoFit() | Out-Null
$wsDe.Columns.Item(16).AutoFit() | Out-Null
